$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = '27.763.74'
$ws.Range("E2").Value = '  -1.76%  '

# Row 3
$ws.Range("D3").Value = '1.896.78'
$ws.Range("E3").Value = '  -1.26%  '

# Row 4
Set-TextValue 'D4' '0.9995'
$ws.Range("E4").Value = '  -0.64%  '

# Row 5
Set-TextValue 'D5' '311.94'
$ws.Range("E5").Value = '  -1.53%  '

# Row 6
Set-TextValue 'D6' '0.9995'
$ws.Range("E6").Value = '  -0.67%  '

# Row 7
Set-TextValue 'D7' '0.4922'
$ws.Range("E7").Value = '  +1.44%  '

# Row 8
Set-TextValue 'D8' '0.3797'
$ws.Range("E8").Value = '  -1.50%  '

# Row 9
Set-TextValue 'D9' '0.07324'
$ws.Range("E9").Value = '  -1.30%  '

# Row 10
Set-TextValue 'D10' '0.9115'
$ws.Range("E10").Value = '  -4.21%  '

# Row 11
Set-TextValue 'D11' '20.65'
$ws.Range("E11").Value = '  -1.49%  '

# Row 12
Set-TextValue 'D12' '0.07625'
$ws.Range("E12").Value = '  -2.36%  '

# Row 13
$ws.Range("D13").Value = '1.903.31'
$ws.Range("E13").Value = '  -1.39%  '

# Row 14
Set-TextValue 'D14' '5.463'
$ws.Range("E14").Value = '  -1.72%  '

# Row 15
Set-TextValue 'D15' '6.654'
$ws.Range("E15").Value = '  -0.13%  '

# Row 16
Set-TextValue 'D16' '91.20'
$ws.Range("E16").Value = '  -1.13%  '

# Row 17
Set-TextValue 'D17' '1.000'
$ws.Range("E17").Value = '  -0.75%  '

# Row 18
Set-TextValue 'D18' '0.000008738'
$ws.Range("E18").Value = '  -1.79%  '

# Row 19
Set-TextValue 'D19' '0.9994'
$ws.Range("E19").Value = '  -0.63%  '

# Row 20
$ws.Range("D20").Value = '27.776.56'
$ws.Range("E20").Value = '  -1.72%  '

# Row 21
Set-TextValue 'D21' '14.48'
$ws.Range("E21").Value = '  -3.91%  '

# Row 22
Set-TextValue 'D22' '5.118'
$ws.Range("E22").Value = '  -1.20%  '

# Row 23
$ws.Range("D23").Value = '2.119.28'
$ws.Range("E23").Value = '  -1.98%  '

# Row 24
$ws.Range("E24").Value = '  -1.71%  '

# Row 25
Set-TextValue 'D25' '153.97'
$ws.Range("E25").Value = '  -1.49%  '

# Row 26
Set-TextValue 'D26' '1.851'
$ws.Range("E26").Value = '  -4.23%  '

# Row 27
Set-TextValue 'D27' '2.176'
$ws.Range("E27").Value = '  +2.69%  '

# Row 28
$ws.Range("E28").Value = '  -1.38%  '

# Row 29
Set-TextValue 'D29' '115.42'
$ws.Range("E29").Value = '  -1.69%  '

# Row 30
Set-TextValue 'D30' '4.887'
$ws.Range("E30").Value = '  -2.88%  '

# Row 31
Set-TextValue 'D31' '0.08935'
$ws.Range("E31").Value = '  +0.23%  '

# Row 32
Set-TextValue 'D32' '3.228'
$ws.Range("E32").Value = '  -4.12%  '

# Row 33
Set-TextValue 'D33' '1.232'
$ws.Range("E33").Value = '  -1.69%  '

# Row 34
Set-TextValue 'D34' '0.7670'
$ws.Range("E34").Value = '  -1.40%  '

# Row 35
Set-TextValue 'D35' '4.643'
$ws.Range("E35").Value = '  -1.01%  '

# Row 36
Set-TextValue 'D36' '2.563'
$ws.Range("E36").Value = '  -8.00%  '

# Row 37
Set-TextValue 'D37' '0.02044'
$ws.Range("E37").Value = '  -0.70%  '

# Row 38
$ws.Range("E38").Value = '  -2.80%  '

# Row 39
$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 'D39' '0.5504'
$ws.Range("E39").Value = '  -1.73%  '

# Row 40
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D40' '0.05284'
$ws.Range("E40").Value = '  -1.81%  '

# Row 41
$ws.Range("E41").Value = '  -1.57%  '

# Row 42
Set-TextValue 'D42' '6.900'
$ws.Range("E42").Value = '  -3.24%  '

# Row 43
Set-TextValue 'D43' '8.574'
$ws.Range("E43").Value = '  -0.22%  '

# Row 44
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D44' '112.82'
$ws.Range("E44").Value = '  +4.74%  '

# Row 45
$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D45' '0.1522'
$ws.Range("E45").Value = '  -1.12%  '

# Row 46
Set-TextValue 'D46' '10.67'
$ws.Range("E46").Value = '  -1.41%  '

# Row 47
Set-TextValue 'D47' '0.4801'
$ws.Range("E47").Value = '  -2.59%  '

# Row 48
Set-TextValue 'D48' '0.9994'
$ws.Range("E48").Value = '  -0.73%  '

# Row 49
Set-TextValue 'D49' '1.635'
$ws.Range("E49").Value = '  -2.72%  '

# Row 50
Set-TextValue 'D50' '67.53'
$ws.Range("E50").Value = '  -2.79%  '

# Row 51
Set-TextValue 'D51' '0.06056'
$ws.Range("E51").Value = '  -1.58%  '
